$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("L4").Value = "not supported"
